# Rename the first four worksheets and delete the last four (now-redundant)
# data worksheets. The last four sheets were exact duplicates of the data
# contained in the first four sheets, and the first four sheets are renamed
# to reflect what the data actually represents (ro_* / code_* placement).

$wb = $excel.ActiveWorkbook

# Rename sheets 1-4
$wb.Worksheets.Item(1).Name = "ro_FLASH-code_FLASH"
$wb.Worksheets.Item(2).Name = "ro_FLASH-code_CCM"
$wb.Worksheets.Item(3).Name = "ro_CCM-code_FLASH"
$wb.Worksheets.Item(4).Name = "ro_CCL-code_CCM"

# Delete the redundant trailing sheets (data_RAM code_FLASH, data_RAM code_CCM,
# data_CCM code_FLASH, data_CCM code_CCM). Delete from the end so indices stay
# valid as we go.
$wb.Worksheets.Item("data_CCM code_CCM").Delete()
$wb.Worksheets.Item("data_CCM code_FLASH").Delete()
$wb.Worksheets.Item("data_RAM code_CCM").Delete()
$wb.Worksheets.Item("data_RAM code_FLASH").Delete()
